$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: becomes the old row-18 record (Grön jordtunga), with a new
#     Taxonsorteringsordning (B) value that differs from a pure swap ---
$ws.Range("A17").Value = 111961716
$ws.Range("B17").Value = 81207
$ws.Range("D17").Value = "LC"
$ws.Range("E17").Value = 5046
$ws.Range("F17").Value = "Grön jordtunga"
$ws.Range("G17").Value = "Microglossum viride"
$ws.Range("H17").Value = "(Pers.:Fr.) Gillet"

# I17 ("Antal") holds the text "2" - force text so it matches the sheet's
# existing text-typed "Antal" column instead of becoming a number.
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = "2"
$ws.Range("I17").NumberFormat = "General"
$ws.Range("I17").Style = "Normal"

$ws.Range("J17").Value = "mycel"
$ws.Range("AF17").Value = "mikroskoperad"
$ws.Range("AX17").Value = "Stefan Phalagorn Bergström, Andreas Estensen, Annika  Carlberg , Ola Elleström, Thomas Strid, Anne Järvinen, Emma Sewell"

# --- Row 18: becomes the old row-17 record (Svartvit taggsvamp), with a
#     new Taxonsorteringsordning (B) value that differs from a pure swap ---
$ws.Range("A18").Value = 111961472
$ws.Range("B18").Value = 90857
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 5448
$ws.Range("F18").Value = "Svartvit taggsvamp"
$ws.Range("G18").Value = "Phellodon connatus"
$ws.Range("H18").Value = "(Schultz) nom.prov"
$ws.Range("I18").Value = ""
$ws.Range("J18").Value = ""
$ws.Range("AF18").Value = ""
$ws.Range("AX18").Value = "Stefan Phalagorn Bergström, Annika  Carlberg , Andreas Estensen, Ola Elleström, Anne Järvinen, Emma Sewell, Thomas Strid"

# --- Row 19: only the Taxonsorteringsordning (B) value changes ---
$ws.Range("B19").Value = 89950
